$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserHome")
$ws.Range("G2").ClearContents()
$ws.Range("G2").Select()
